$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure the NMR_METABOLITE_DATA block (rows 76-84 -> 76-86) ---
# 1) Remove the old "Metabolites" value row (old row 78); this shifts rows 79+ up by one.
$ws.Rows.Item(78).Delete()

# 2) Insert two fresh blank rows where the new "#.collate"/"Metabolites" pair
#    (with the new *#.optional_headers column) needs to go - right before the
#    old "*#.fields" header block (which is now at row 82 after the delete above).
$ws.Rows.Item(82).Insert()
$ws.Rows.Item(82).Insert()

# --- Populate the final contents of the restructured block ---

# Row 76: header row for #.collate (unchanged layout/columns)
$ws.Range("A76").Value = "#tags"
$ws.Range("B76").Value = "#NMR_METABOLITE_DATA.id"
$ws.Range("C76").Value = "#.collate"
$ws.Range("D76").Value = "*#.headers"
$ws.Range("E76").Value = "#.required"
$ws.Range("F76").Value = "*#.sort_by"
$ws.Range("G76").Value = "#.sort_order"
$ws.Range("H76").Value = "#.table"
$ws.Range("I76").Value = "#.value_type"
$ws.Range("J76").Value = "#.values_to_str"

# Row 77: "Data" values - resonance_assignment renamed to assignment
$ws.Range("B77").Value = "Data"
$ws.Range("C77").Value = "assignment"
$ws.Range("D77").Value = "`"Metabolite`"=assignment,entity.id=intensity"
$ws.Range("E77").Value = "'True"
$ws.Range("F77").Value = "assignment"
$ws.Range("G77").Value = "ascending"
$ws.Range("H77").Value = "measurement"
$ws.Range("I77").Value = "matrix"
$ws.Range("J77").Value = "'True"

# Row 78 is blank (separator)

# Row 79: header row for *#.exclusion_headers / #.fields_to_headers
$ws.Range("A79").Value = "#tags"
$ws.Range("B79").Value = "#NMR_METABOLITE_DATA.id"
$ws.Range("C79").Value = "*#.exclusion_headers"
$ws.Range("D79").Value = "#.fields_to_headers"
$ws.Range("E79").Value = "*#.headers"
$ws.Range("F79").Value = "#.required"
$ws.Range("G79").Value = "*#.sort_by"
$ws.Range("H79").Value = "#.sort_order"
$ws.Range("I79").Value = "#.table"
$ws.Range("J79").Value = "#.value_type"
$ws.Range("K79").Value = "#.values_to_str"

# Row 80: "Extended" values - new assignment/assignment%method fields added
$ws.Range("B80").Value = "Extended"
$ws.Range("C80").Value = "id,intensity,intensity%type,intensity%units,assignment,assignment%method,entity.id,protocol.id,base_inchi,representative_inchi,isotopic_inchi,peak_description,peak_pattern,proton_count,transient_peak,transient_peak%type"
$ws.Range("D80").Value = "'True"
$ws.Range("E80").Value = "`"Metabolite`"=assignment,`"sample_id`"=entity.id"
$ws.Range("F80").Value = "'False"
$ws.Range("G80").Value = "assignment"
$ws.Range("H80").Value = "ascending"
$ws.Range("I80").Value = "measurement"
$ws.Range("J80").Value = "matrix"
$ws.Range("K80").Value = "'True"

# Row 81 is blank (separator)

# Row 82: new header row for #.collate / *#.optional_headers ("Metabolites" group)
$ws.Range("A82").Value = "#tags"
$ws.Range("B82").Value = "#NMR_METABOLITE_DATA.id"
$ws.Range("C82").Value = "#.collate"
$ws.Range("D82").Value = "*#.headers"
$ws.Range("E82").Value = "*#.optional_headers"
$ws.Range("F82").Value = "#.required"
$ws.Range("G82").Value = "*#.sort_by"
$ws.Range("H82").Value = "#.sort_order"
$ws.Range("I82").Value = "#.table"
$ws.Range("J82").Value = "#.value_type"
$ws.Range("K82").Value = "#.values_to_str"

# Row 83: "Metabolites" values (moved here with new optional_headers column)
$ws.Range("B83").Value = "Metabolites"
$ws.Range("C83").Value = "assignment"
$ws.Range("D83").Value = "`"Metabolite`"=assignment"
$ws.Range("E83").Value = "assignment%method,base_inchi,representative_inchi,isotopic_inchi,peak_description,peak_pattern,proton_count,transient_peak,transient_peak%type"
$ws.Range("F83").Value = "'True"
$ws.Range("G83").Value = "assignment"
$ws.Range("H83").Value = "ascending"
$ws.Range("I83").Value = "measurement"
$ws.Range("J83").Value = "matrix"
$ws.Range("K83").Value = "'True"

# Row 84 is blank (separator) - clear the stale "*#.fields" header that
# landed here after the row shift from the delete/insert above.
$ws.Rows.Item(84).ClearContents()

# Row 85: header row for *#.fields (unchanged)
$ws.Range("A85").Value = "#tags"
$ws.Range("B85").Value = "#NMR_METABOLITE_DATA.id"
$ws.Range("C85").Value = "*#.fields"
$ws.Range("D85").Value = "#.required"
$ws.Range("E85").Value = "#.table"
$ws.Range("F85").Value = "#.value_type"

# Row 86: "Units" values (unchanged)
$ws.Range("B86").Value = "Units"
$ws.Range("C86").Value = "intensity%type"
$ws.Range("D86").Value = "'True"
$ws.Range("E86").Value = "measurement"
$ws.Range("F86").Value = "str"
